$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.979.52'
$ws.Range('E2').Value = '  +1.86%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.307.36'
$ws.Range('E3').Value = '  +1.64%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.70'
$ws.Range('E5').Value = '  +1.92%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.14'
$ws.Range('E6').Value = '  +5.27%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.505'
$ws.Range('E7').Value = '  +2.43%  '

$ws.Range('E9').Value = '  +4.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.88'
$ws.Range('E10').Value = '  +4.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0796'
$ws.Range('E11').Value = '  +0.84%  '

$ws.Range('E12').Value = '  +4.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.02'
$ws.Range('E13').Value = '  +16.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.86'
$ws.Range('E14').Value = '  +3.07%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.687.66'
$ws.Range('E15').Value = '  +2.44%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.341.86'
$ws.Range('E16').Value = '  +2.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.816'
$ws.Range('E17').Value = '  +4.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.925.20'
$ws.Range('E18').Value = '  +1.88%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.46'
$ws.Range('E19').Value = '  +6.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.14'
$ws.Range('E20').Value = '  +2.48%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  +1.39%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.68'
$ws.Range('E22').Value = '  +1.75%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.22'
$ws.Range('E23').Value = '  +1.38%  '

$ws.Range('E24').Value = '  +12.60%  '

$ws.Range('E25').Value = '  +0.50%  '

$ws.Range('E26').Value = '  -0.44%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.77'
$ws.Range('E27').Value = '  +3.65%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -4.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.37'
$ws.Range('E29').Value = '  -0.33%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.00'
$ws.Range('E30').Value = '  -0.49%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.15'
$ws.Range('E31').Value = '  +0.59%  '

$ws.Range('E32').Value = '  +0.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.00'
$ws.Range('E33').Value = '  +2.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.63'
$ws.Range('E34').Value = '  +1.88%  '

$ws.Range('E35').Value = '  +4.05%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.97'
$ws.Range('E36').Value = '  +2.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0690'
$ws.Range('E37').Value = '  +0.64%  '

$ws.Range('E38').Value = '  +3.54%  '

$ws.Range('E39').Value = '  +1.28%  '

$ws.Range('E40').Value = '  +3.57%  '

$ws.Range('E41').Value = '  +0.62%  '

$ws.Range('E42').Value = '  -6.11%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.002.40'
$ws.Range('E43').Value = '  +2.10%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0285'
$ws.Range('E44').Value = '  +2.77%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.21'
$ws.Range('E45').Value = '  +6.92%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.34'
$ws.Range('E46').Value = '  -0.71%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.83'
$ws.Range('E47').Value = '  +1.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.49'
$ws.Range('E48').Value = '  +6.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.530.38'
$ws.Range('E49').Value = '  +1.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').Value = '  +4.26%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.56'
$ws.Range('E51').Value = '  +0.89%  '
